$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value (keeps the "looks like a number" string as TEXT,
# i.e. t="s" in the OOXML) into $cellRef without leaving a numberformat/style
# behind on the target cell itself. We stage the value (with a leading
# apostrophe so Excel treats it as text) in a scratch cell, copy it, paste
# only the (unformatted) value into the destination, then wipe the scratch
# cell completely.
function Set-TextValue {
    param([string]$cellRef, [string]$text)
    $scratch = $ws.Range("ZZ1")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---- Row 2 : existing row, only E2 changes (number -> text "MID0001") ----
# (written later, in the same relative order the shared strings were first
# introduced - see below)

# ---- Row 3 (new) ----
Set-TextValue "A3" "1"
Set-TextValue "C3" "1"
Set-TextValue "D3" "3"
$ws.Range("G3").Value = "test"
$ws.Range("K3").Value = "RP-8/1/2020"

# ---- Row 4 (new) ----
Set-TextValue "A4" "1"
Set-TextValue "C4" "3"
$ws.Range("G4").Value = "test"
$ws.Range("K4").Value = "RP-8/1/2020"

# ---- Row 5 (new) ----
$ws.Range("G5").Value = "test"
$ws.Range("K5").Value = "RP-8/1/2020"

# ---- Row 6 (new) ----
$ws.Range("G6").Value = "test"
$ws.Range("K6").Value = "RP-8/1/2020"

# ---- Row 7 (new) ----
$ws.Range("G7").Value = "test2"
$ws.Range("K7").Value = "RP-8/1/2020"

# ---- Med_ID text column (E) ----
$ws.Range("E3").Value = "MID0004"
$ws.Range("E9").Value = "MID0004"
$ws.Range("E4").Value = "MID0005"
$ws.Range("E7").Value = "MID0003"
$ws.Range("E8").Value = "MID0003"
$ws.Range("E2").Value = "MID0001"
$ws.Range("E5").Value = "MID0002"
$ws.Range("E6").Value = "MID0002"

$ws.Range("G8").Value = "test"
$ws.Range("G9").Value = "test"

# ---- Numeric MRNumber / Day_Type columns for rows 5-9 (plain numbers) ----
$ws.Range("A5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("A6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("A8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("A9").Value = 1
$ws.Range("C9").Value = 1

# ---- Reason_For_Change_Med (F) for rows 3-4 only ----
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0

# ---- Daily_Med_Dose_Mg (H) ----
$ws.Range("H3").Value = 1.5
$ws.Range("H4").Value = 3.1
$ws.Range("H5").Value = 2.1
$ws.Range("H6").Value = 1.3
$ws.Range("H7").Value = 3.2
$ws.Range("H8").Value = 5
$ws.Range("H9").Value = 3

# ---- Date column (B) : same date style as B2, reuse it via copy/paste ----
$ws.Range("B2").Copy()
$ws.Range("B3:B9").PasteSpecial(-4122)
$ws.Range("B3").Value = 43893
$ws.Range("B4").Value = 43952
$ws.Range("B5").Value = 43984
$ws.Range("B6").Value = 44014
$ws.Range("B7").Value = 43863
$ws.Range("B8").Value = 43893
$ws.Range("B9").Value = 43894

[void]$ws.Range("H9").Select()
